$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain stored as text,
# since some new values (e.g. "0.9997", "15.02") would otherwise be
# auto-coerced by Excel into numbers. We flip NumberFormat to Text,
# assign, then restore the default "Normal" style so no visible/
# persisted formatting change is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '25.832.37'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '1.736.12'
$ws.Range('E3').Value = '  -1.27%  '
Set-TextValue $ws.Range('D4') '0.9997'
Set-TextValue $ws.Range('D5') '230.42'
$ws.Range('E5').Value = '  -2.82%  '
Set-TextValue $ws.Range('D6') '0.9996'
$ws.Range('E6').Value = '  -0.06%  '
Set-TextValue $ws.Range('D7') '0.5220'
$ws.Range('E7').Value = '  -0.06%  '
Set-TextValue $ws.Range('D8') '0.2752'
$ws.Range('E8').Value = '  +1.67%  '
Set-TextValue $ws.Range('D9') '39.31'
$ws.Range('E9').Value = '  -2.89%  '
Set-TextValue $ws.Range('D10') '0.06140'
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').Value = '1.734.71'
$ws.Range('E11').Value = '  -1.98%  '
Set-TextValue $ws.Range('D12') '0.07032'
$ws.Range('E12').Value = '  +0.12%  '
Set-TextValue $ws.Range('D13') '15.02'
$ws.Range('E13').Value = '  -4.36%  '
Set-TextValue $ws.Range('D14') '0.6351'
$ws.Range('E14').Value = '  -3.32%  '
Set-TextValue $ws.Range('D15') '4.527'
$ws.Range('E15').Value = '  +0.97%  '
Set-TextValue $ws.Range('D16') '76.71'
$ws.Range('E16').Value = '  -1.86%  '
Set-TextValue $ws.Range('D17') '0.9996'
$ws.Range('E17').Value = '  -0.04%  '
Set-TextValue $ws.Range('D18') '0.9995'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '25.815.67'
$ws.Range('E19').Value = '  -0.71%  '
Set-TextValue $ws.Range('D20') '11.46'
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = '1.955.20'
$ws.Range('E22').Value = '  -1.44%  '
Set-TextValue $ws.Range('D23') '4.185'
$ws.Range('E23').Value = '  +2.27%  '
Set-TextValue $ws.Range('D24') '8.747'
$ws.Range('E24').Value = '  +4.06%  '
Set-TextValue $ws.Range('D25') '5.155'
$ws.Range('E25').Value = '  -0.43%  '
Set-TextValue $ws.Range('D26') '139.24'
$ws.Range('E26').Value = '  +1.44%  '
Set-TextValue $ws.Range('D27') '1.502'
$ws.Range('E27').Value = '  +1.37%  '
Set-TextValue $ws.Range('D28') '15.01'
$ws.Range('E28').Value = '  -0.86%  '
$ws.Range('E29').Value = '  -2.81%  '
Set-TextValue $ws.Range('D30') '101.92'
$ws.Range('E30').Value = '  -0.61%  '
Set-TextValue $ws.Range('D31') '0.08300'
$ws.Range('E31').Value = '  -1.31%  '
Set-TextValue $ws.Range('D32') '3.707'
$ws.Range('E32').Value = '  +0.22%  '
Set-TextValue $ws.Range('D33') '3.496'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('E34').Value = '  +0.96%  '
Set-TextValue $ws.Range('D35') '2.603'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('E36').Value = '  -2.42%  '
Set-TextValue $ws.Range('D37') '0.6151'
$ws.Range('E37').Value = '  +1.22%  '
Set-TextValue $ws.Range('D38') '2.669'
$ws.Range('E38').Value = '  -2.45%  '
Set-TextValue $ws.Range('D39') '0.01569'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range('D40') '0.9990'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D41') '1.903'
$ws.Range('E41').Value = '  -2.36%  '
Set-TextValue $ws.Range('D42') '100.01'
$ws.Range('E42').Value = '  -2.82%  '
$ws.Range('E43').Value = '  -1.25%  '
Set-TextValue $ws.Range('D44') '5.005'
$ws.Range('E44').Value = '  +1.46%  '
Set-TextValue $ws.Range('D45') '0.7202'
$ws.Range('E45').Value = '  -4.34%  '
$ws.Range('E46').Value = '  -2.82%  '
Set-TextValue $ws.Range('D47') '0.1127'
$ws.Range('E47').Value = '  +0.77%  '
Set-TextValue $ws.Range('D48') '6.162'
$ws.Range('E48').Value = '  +0.99%  '
Set-TextValue $ws.Range('D49') '53.12'
$ws.Range('E49').Value = '  +0.91%  '
Set-TextValue $ws.Range('D50') '29.97'
$ws.Range('E50').Value = '  -0.56%  '
Set-TextValue $ws.Range('D51') '7.615'
$ws.Range('E51').Value = '  +2.48%  '
